$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row fixes: "first name" / "last name" -> "first_name" / "last_name" ---
$ws.Cells.Item(1, 3).Value = "first_name"
$ws.Cells.Item(1, 4).Value = "last_name"

# --- Existing rows 2-11: columns C (first_name) & D (last_name) now hold "userN" instead of "usernameN" ---
$userNames = @("user1","user2","user3","user4","user5","user6","user7","user8","user9","user10")
for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $userNames[$i]
    $ws.Cells.Item($row, 4).Value = $userNames[$i]
}

# --- New rows 12-21: username11-20 / password11-20 / user11-20 / user11-20 / email (cycled) / user=true / manager=false ---
for ($i = 1; $i -le 10; $i++) {
    $row = 11 + $i
    $n = 10 + $i
    $ws.Cells.Item($row, 1).Value = "username$n"
    $ws.Cells.Item($row, 2).Value = "password$n"
    $ws.Cells.Item($row, 3).Value = "user$n"
    $ws.Cells.Item($row, 4).Value = "user$n"
    $ws.Cells.Item($row, 6).Value = $true
    $ws.Cells.Item($row, 7).Value = $false
}

# Hyperlinks for the new rows' email column (E), added in the same order the
# original workbook's diff lists them so relationship ids line up (rId11..rId20)
$hyperlinkOrder = @(
    @{Row=12; Email="username1@gmail.com"},
    @{Row=13; Email="username2@gmail.com"},
    @{Row=14; Email="username3@gmail.com"},
    @{Row=15; Email="username4@gmail.com"},
    @{Row=18; Email="username7@gmail.com"},
    @{Row=21; Email="username10@gmail.com"},
    @{Row=16; Email="username5@gmail.com"},
    @{Row=19; Email="username8@gmail.com"},
    @{Row=17; Email="username6@gmail.com"},
    @{Row=20; Email="username9@gmail.com"}
)

foreach ($link in $hyperlinkOrder) {
    $cell = $ws.Cells.Item($link.Row, 5)
    # Pre-set the display text so Hyperlinks.Add reuses the existing shared
    # string instead of minting a new "mailto:..." one.
    $cell.Value = $link.Email
    $ws.Hyperlinks.Add($cell, "mailto:" + $link.Email) | Out-Null
    # Match the style used by the existing email hyperlink cells (E2:E11)
    $cell.Style = $ws.Cells.Item(2, 5).Style
}

# --- Selection ---
$ws.Range("I19").Select() | Out-Null
